# Auto-generated Excel COM-interop script to apply the edit described by the diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet 1: "展览" (Exhibitions) -----
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 1).Value = 1
$ws1.Cells.Item(2, 2).Value = "2024-03-16"
$ws1.Cells.Item(2, 3).Value = "杭州·ComicMe动漫嘉年华"
$ws1.Cells.Item(2, 4).Value = "长江南路336号 白马湖国际会展中心"
$ws1.Cells.Item(2, 5).Value = "2024.03.16 09:00-03.17 17:00"
$ws1.Cells.Item(2, 6).Value = 2592
$ws1.Cells.Item(2, 7).Value = 60
$ws1.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81214"
$ws1.Cells.Item(2, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/iN5DafVM1705907838033.jpeg"

$ws1.Cells.Item(3, 1).Value = 2
$ws1.Cells.Item(3, 2).Value = "2024-03-16"
$ws1.Cells.Item(3, 3).Value = "杭州·ComicMe动漫嘉年华 · 胡良伟专场"
$ws1.Cells.Item(3, 4).Value = "长江南路336号 白马湖国际会展中心"
$ws1.Cells.Item(3, 5).Value = "2024.03.16 12:30-03.16 16:30"
$ws1.Cells.Item(3, 6).Value = 340
$ws1.Cells.Item(3, 7).Value = "已售罄"
$ws1.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81382"
$ws1.Cells.Item(3, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/4UJgRWeM1706151833915.jpeg"

$ws1.Cells.Item(4, 1).Value = 3
$ws1.Cells.Item(4, 2).Value = "2024-03-16"
$ws1.Cells.Item(4, 3).Value = "杭州·OZ·富坚义博only"
$ws1.Cells.Item(4, 4).Value = "北干街道萧杭路689号浙农东巢艺术公园 Fashion Bund时尚外滩艺术中心"
$ws1.Cells.Item(4, 5).Value = "2024.03.16 10:00-03.16 17:00"
$ws1.Cells.Item(4, 6).Value = 346
$ws1.Cells.Item(4, 7).Value = 88
$ws1.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81151"
$ws1.Cells.Item(4, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/CxqdeAPa1705658329588.jpeg"

$ws1.Cells.Item(5, 1).Value = 4
$ws1.Cells.Item(5, 2).Value = "2024-03-16"
$ws1.Cells.Item(5, 3).Value = "杭州·SST动漫嘉年华"
$ws1.Cells.Item(5, 4).Value = "沈半路171号 Tcar汽车文化主题公园"
$ws1.Cells.Item(5, 5).Value = "2024.03.16 09:00-03.17 17:00"
$ws1.Cells.Item(5, 6).Value = 1449
$ws1.Cells.Item(5, 7).Value = 68
$ws1.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81196"
$ws1.Cells.Item(5, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/yFyT2uAT1705897787652.jpeg"

$ws1.Cells.Item(6, 1).Value = 5
$ws1.Cells.Item(6, 2).Value = "2024-03-16"
$ws1.Cells.Item(6, 3).Value = "杭州·排球少年*蓝锁ONLY"
$ws1.Cells.Item(6, 4).Value = "亚太路湘湖3期东南侧约290米 原创壹号羽毛球馆"
$ws1.Cells.Item(6, 5).Value = "2024.03.16 10:00-03.16 17:00"
$ws1.Cells.Item(6, 6).Value = 1124
$ws1.Cells.Item(6, 7).Value = 60
$ws1.Cells.Item(6, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81075"
$ws1.Cells.Item(6, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/9AL6kYuj1705634962275.jpeg"

$ws1.Cells.Item(7, 1).Value = 6
$ws1.Cells.Item(7, 2).Value = "2024-03-16"
$ws1.Cells.Item(7, 3).Value = "杭州·春和景明代号鸢only"
$ws1.Cells.Item(7, 4).Value = "金沙大道681号 金沙湖大剧院"
$ws1.Cells.Item(7, 5).Value = "2024.03.16 09:30-03.16 16:00"
$ws1.Cells.Item(7, 6).Value = 324
$ws1.Cells.Item(7, 7).Value = "已售罄"
$ws1.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81894"
$ws1.Cells.Item(7, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/1RX6wnFN1708323470900.png"

$ws1.Cells.Item(8, 1).Value = 7
$ws1.Cells.Item(8, 2).Value = "2024-03-16"
$ws1.Cells.Item(8, 3).Value = "杭州·百鬼夜行·咒术回战only"
$ws1.Cells.Item(8, 4).Value = "长生路18号 梅地亚宾馆"
$ws1.Cells.Item(8, 5).Value = "2024.03.16 09:00-03.16 17:00"
$ws1.Cells.Item(8, 6).Value = 531
$ws1.Cells.Item(8, 7).Value = 79
$ws1.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81478"
$ws1.Cells.Item(8, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/4weHdCdk1706495040356.jpeg"

$ws1.Cells.Item(9, 1).Value = 8
$ws1.Cells.Item(9, 2).Value = "2024-03-16"
$ws1.Cells.Item(9, 3).Value = "杭州·第十届次元鹿角动漫游戏展（取消）"
$ws1.Cells.Item(9, 4).Value = "万融城3幢1楼 头号玩家数字运动俱乐部"
$ws1.Cells.Item(9, 5).Value = "2024.03.16 10:00-03.17 17:00"
$ws1.Cells.Item(9, 6).Value = 1163
$ws1.Cells.Item(9, 7).Value = "不可售"
$ws1.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81136"
$ws1.Cells.Item(9, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/w8iHjfOA1705651976885.jpeg"

$ws1.Cells.Item(10, 1).Value = 9
$ws1.Cells.Item(10, 2).Value = "2024-03-16"
$ws1.Cells.Item(10, 3).Value = "杭州·筑梦城堡巡回展降临之章（取消）"
$ws1.Cells.Item(10, 4).Value = "大岭山路156号 爱丽芬城堡"
$ws1.Cells.Item(10, 5).Value = "2024.03.16 10:00-03.16 17:00"
$ws1.Cells.Item(10, 6).Value = 16
$ws1.Cells.Item(10, 7).Value = "不可售"
$ws1.Cells.Item(10, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81217"
$ws1.Cells.Item(10, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/oEILAgir1705908585941.png"

$ws1.Cells.Item(11, 1).Value = 10
$ws1.Cells.Item(11, 2).Value = "2024-03-16"
$ws1.Cells.Item(11, 3).Value = "杭州·造梦探险家——次元茶话会"
$ws1.Cells.Item(11, 4).Value = "临平街道北沙西路156-1号 杭州临平遇上设计师酒店"
$ws1.Cells.Item(11, 5).Value = "2024.03.16 10:00-03.16 17:00"
$ws1.Cells.Item(11, 6).Value = 111
$ws1.Cells.Item(11, 7).Value = 38
$ws1.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81914"
$ws1.Cells.Item(11, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/HHHVcvcC1709015213282.png"

$ws1.Cells.Item(12, 1).Value = 11
$ws1.Cells.Item(12, 2).Value = "2024-03-17"
$ws1.Cells.Item(12, 3).Value = "杭州·ComicMe动漫嘉年华 · 马正阳专场"
$ws1.Cells.Item(12, 4).Value = "长江南路336号 白马湖国际会展中心"
$ws1.Cells.Item(12, 5).Value = "2024.03.17 12:30-03.17 16:30"
$ws1.Cells.Item(12, 6).Value = 549
$ws1.Cells.Item(12, 7).Value = "已售罄"
$ws1.Cells.Item(12, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81375"
$ws1.Cells.Item(12, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/rNzSh0nq1706149891822.jpeg"

$ws1.Cells.Item(13, 1).Value = 12
$ws1.Cells.Item(13, 2).Value = "2024-03-23"
$ws1.Cells.Item(13, 3).Value = "杭州·AD02动漫展"
$ws1.Cells.Item(13, 4).Value = "浙江省杭州市萧山区奔竞大道353号 国际博览中心"
$ws1.Cells.Item(13, 5).Value = "2024.03.23 10:00-03.24 17:00"
$ws1.Cells.Item(13, 6).Value = 8934
$ws1.Cells.Item(13, 7).Value = 75
$ws1.Cells.Item(13, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80905"
$ws1.Cells.Item(13, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/D3QaPamg1705397424553.jpeg"

$ws1.Cells.Item(14, 1).Value = 13
$ws1.Cells.Item(14, 2).Value = "2024-03-23"
$ws1.Cells.Item(14, 3).Value = "杭州·AD02动漫展  青柳尊哉内场票"
$ws1.Cells.Item(14, 4).Value = "浙江省杭州市萧山区奔竞大道353号 国际博览中心"
$ws1.Cells.Item(14, 5).Value = "2024.03.23 10:00-03.23 17:00"
$ws1.Cells.Item(14, 6).Value = 386
$ws1.Cells.Item(14, 7).Value = 528
$ws1.Cells.Item(14, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81503"
$ws1.Cells.Item(14, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/OmqxboDC1706522627528.jpeg"

$ws1.Cells.Item(15, 1).Value = 14
$ws1.Cells.Item(15, 2).Value = "2024-03-23"
$ws1.Cells.Item(15, 3).Value = "杭州·AD02动漫展--卡琳娜签售票"
$ws1.Cells.Item(15, 4).Value = "浙江省杭州市萧山区奔竞大道353号 国际博览中心"
$ws1.Cells.Item(15, 5).Value = "2024.03.23 09:30-03.23 17:00"
$ws1.Cells.Item(15, 6).Value = 2497
$ws1.Cells.Item(15, 7).Value = "已售罄"
$ws1.Cells.Item(15, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81941"
$ws1.Cells.Item(15, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/CZjxY9ZC1708416661613.jpeg"

$ws1.Cells.Item(16, 1).Value = 15
$ws1.Cells.Item(16, 2).Value = "2024-03-24"
$ws1.Cells.Item(16, 3).Value = "杭州·AD02动漫展  岩永彻也内场票"
$ws1.Cells.Item(16, 4).Value = "浙江省杭州市萧山区奔竞大道353号 国际博览中心"
$ws1.Cells.Item(16, 5).Value = "2024.03.24 10:00-03.24 17:00"
$ws1.Cells.Item(16, 6).Value = 253
$ws1.Cells.Item(16, 7).Value = 528
$ws1.Cells.Item(16, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81239"
$ws1.Cells.Item(16, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/hww9WUpD1705914756383.jpeg"

$ws1.Cells.Item(17, 1).Value = 16
$ws1.Cells.Item(17, 2).Value = "2024-03-24"
$ws1.Cells.Item(17, 3).Value = "杭州·AD02动漫展--亦之紫F、L句号内场票"
$ws1.Cells.Item(17, 4).Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws1.Cells.Item(17, 5).Value = "2024.03.24 12:00-03.24 16:00"
$ws1.Cells.Item(17, 6).Value = 177
$ws1.Cells.Item(17, 7).Value = 258
$ws1.Cells.Item(17, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81836"
$ws1.Cells.Item(17, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/ecrRfQce1707375167618.jpeg"

$ws1.Cells.Item(18, 1).Value = 17
$ws1.Cells.Item(18, 2).Value = "2024-03-24"
$ws1.Cells.Item(18, 3).Value = "杭州·AD02动漫展--钟晨瑶内场票"
$ws1.Cells.Item(18, 4).Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws1.Cells.Item(18, 5).Value = "2024.03.24 09:30-03.24 17:00"
$ws1.Cells.Item(18, 6).Value = 469
$ws1.Cells.Item(18, 7).Value = "已售罄"
$ws1.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81820"
$ws1.Cells.Item(18, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/aHRmCxr31707296105225.jpeg"

$ws1.Cells.Item(19, 1).Value = 18
$ws1.Cells.Item(19, 2).Value = "2024-03-30"
$ws1.Cells.Item(19, 3).Value = "杭州·Look Look动漫嘉年华"
$ws1.Cells.Item(19, 4).Value = "聚业路27号电魂大厦B座1楼 电魂自在里文化空间(硅谷书房)"
$ws1.Cells.Item(19, 5).Value = "2024.03.30 10:00-03.31 17:30"
$ws1.Cells.Item(19, 6).Value = 609
$ws1.Cells.Item(19, 7).Value = 52.2
$ws1.Cells.Item(19, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81757"
$ws1.Cells.Item(19, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/anrpU73c1707106069934.jpeg"

$ws1.Cells.Item(20, 1).Value = 19
$ws1.Cells.Item(20, 2).Value = "2024-03-30"
$ws1.Cells.Item(20, 3).Value = "杭州·二次元拾梦漫展（取消）"
$ws1.Cells.Item(20, 4).Value = "转塘街道创意路1号 艺创小镇凤凰创意大厦"
$ws1.Cells.Item(20, 5).Value = "2024.03.30 10:00-03.31 18:00"
$ws1.Cells.Item(20, 6).Value = 71
$ws1.Cells.Item(20, 7).Value = "不可售"
$ws1.Cells.Item(20, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81902"
$ws1.Cells.Item(20, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/OfonGbvl1708325977132.jpeg"

$ws1.Cells.Item(21, 1).Value = 20
$ws1.Cells.Item(21, 2).Value = "2024-03-30"
$ws1.Cells.Item(21, 3).Value = "杭州·幻想物语动漫游戏展"
$ws1.Cells.Item(21, 4).Value = "富春路80号(甬江路地铁站A口旁) 杭州全民健身中心"
$ws1.Cells.Item(21, 5).Value = "2024.03.30 10:00-03.31 17:00"
$ws1.Cells.Item(21, 6).Value = 1166
$ws1.Cells.Item(21, 7).Value = 99
$ws1.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81694"
$ws1.Cells.Item(21, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/y2UW01sH1708324524472.jpeg"

$ws1.Cells.Item(22, 1).Value = 21
$ws1.Cells.Item(22, 2).Value = "2024-03-30"
$ws1.Cells.Item(22, 3).Value = "杭州·排球少年only·春日校庆"
$ws1.Cells.Item(22, 4).Value = "之江路149号 云栖培训基地"
$ws1.Cells.Item(22, 5).Value = "2024.03.30 10:00-03.31 17:00"
$ws1.Cells.Item(22, 6).Value = 1001
$ws1.Cells.Item(22, 7).Value = "已售罄"
$ws1.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81511"
$ws1.Cells.Item(22, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/RDI807mS1708410823039.jpeg"

$ws1.Cells.Item(23, 1).Value = 22
$ws1.Cells.Item(23, 2).Value = "2024-04-04"
$ws1.Cells.Item(23, 3).Value = "杭州·2024ESCC游戏电竞博览会暨新次元微光青春动漫交流会"
$ws1.Cells.Item(23, 4).Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws1.Cells.Item(23, 5).Value = "2024.04.04 09:30-04.05 16:30"
$ws1.Cells.Item(23, 6).Value = 2058
$ws1.Cells.Item(23, 7).Value = "不可售"
$ws1.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81450"
$ws1.Cells.Item(23, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/OfpkJ50P1709548942017.png"

$ws1.Cells.Item(24, 1).Value = 23
$ws1.Cells.Item(24, 2).Value = "2024-04-04"
$ws1.Cells.Item(24, 3).Value = "杭州·ELECTRIC COMIC动漫游戏展"
$ws1.Cells.Item(24, 4).Value = "沈半路171号 T-Car杭州汽车文化主题公园"
$ws1.Cells.Item(24, 5).Value = "2024.04.04 10:00-04.05 17:00"
$ws1.Cells.Item(24, 6).Value = 2135
$ws1.Cells.Item(24, 7).Value = 70
$ws1.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82270"
$ws1.Cells.Item(24, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/Be5zFgv11709543746638.jpeg"

$ws1.Cells.Item(25, 1).Value = 24
$ws1.Cells.Item(25, 2).Value = "2024-04-04"
$ws1.Cells.Item(25, 3).Value = "杭州·创造力动漫游戏嘉年华1.0"
$ws1.Cells.Item(25, 4).Value = "沈半路171号 T-Car杭州汽车文化主题公园"
$ws1.Cells.Item(25, 5).Value = "2024.04.04 10:00-04.05 17:00"
$ws1.Cells.Item(25, 6).Value = 62
$ws1.Cells.Item(25, 7).Value = "不可售"
$ws1.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81078"
$ws1.Cells.Item(25, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/o4cl1vwE1705635692432.jpeg"

$ws1.Cells.Item(26, 1).Value = 25
$ws1.Cells.Item(26, 2).Value = "2024-04-04"
$ws1.Cells.Item(26, 3).Value = "杭州·梦漫星河动漫展"
$ws1.Cells.Item(26, 4).Value = "德胜东路2539号 梦马汽车小镇"
$ws1.Cells.Item(26, 5).Value = "2024.04.04 10:00-04.05 17:00"
$ws1.Cells.Item(26, 6).Value = 1839
$ws1.Cells.Item(26, 7).Value = 58.5
$ws1.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81699"
$ws1.Cells.Item(26, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/sZfZd47Y1706868453434.jpeg"

$ws1.Cells.Item(27, 1).Value = 26
$ws1.Cells.Item(27, 2).Value = "2024-04-04"
$ws1.Cells.Item(27, 3).Value = "杭州·第九届萌次元动漫嘉年华"
$ws1.Cells.Item(27, 4).Value = "长乐路29号五组2幢 杭州运河文化发布中心"
$ws1.Cells.Item(27, 5).Value = "2024.04.04 10:00-04.05 17:00"
$ws1.Cells.Item(27, 6).Value = 240
$ws1.Cells.Item(27, 7).Value = "不可售"
$ws1.Cells.Item(27, 8).Value = "https://show.bilibili.com/platform/detail.html?id=78866"
$ws1.Cells.Item(27, 9).Value = "//i1.hdslb.com/bfs/openplatform/202311/8jSeAOZH1700636327971.jpeg"

$ws1.Cells.Item(28, 1).Value = 27
$ws1.Cells.Item(28, 2).Value = "2024-04-05"
$ws1.Cells.Item(28, 3).Value = "杭州·ESCC电竞博览会 倒霉死勒内场票"
$ws1.Cells.Item(28, 4).Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws1.Cells.Item(28, 5).Value = "2024.04.05 09:30-04.05 16:30"
$ws1.Cells.Item(28, 6).Value = 1921
$ws1.Cells.Item(28, 7).Value = "已售罄"
$ws1.Cells.Item(28, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81681"
$ws1.Cells.Item(28, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/suXI547M1706862164353.png"

$ws1.Cells.Item(29, 1).Value = 28
$ws1.Cells.Item(29, 2).Value = "2024-04-05"
$ws1.Cells.Item(29, 3).Value = "杭州·ESCC电竞博览会·钱琛签售礼包"
$ws1.Cells.Item(29, 4).Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws1.Cells.Item(29, 5).Value = "2024.04.05 09:30-04.05 16:30"
$ws1.Cells.Item(29, 6).Value = 480
$ws1.Cells.Item(29, 7).Value = "已售罄"
$ws1.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81680"
$ws1.Cells.Item(29, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/noqtqw701706861615316.png"

$ws1.Cells.Item(30, 1).Value = 29
$ws1.Cells.Item(30, 2).Value = "2024-04-05"
$ws1.Cells.Item(30, 3).Value = "杭州·第36届 中二病 原神x星穹only"
$ws1.Cells.Item(30, 4).Value = "康候圣街99号 顺丰创新中心"
$ws1.Cells.Item(30, 5).Value = "2024.04.05 11:00-04.06 17:00"
$ws1.Cells.Item(30, 6).Value = 534
$ws1.Cells.Item(30, 7).Value = 60
$ws1.Cells.Item(30, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82313"
$ws1.Cells.Item(30, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/3j66NHdg1709619696758.jpeg"

$ws1.Cells.Item(31, 1).Value = 30
$ws1.Cells.Item(31, 2).Value = "2024-04-13"
$ws1.Cells.Item(31, 3).Value = "杭州·ACG发色only-黑白两色"
$ws1.Cells.Item(31, 4).Value = "康候圣街99号 顺丰创新中心"
$ws1.Cells.Item(31, 5).Value = "2024.04.13 09:00-04.13 18:00"
$ws1.Cells.Item(31, 6).Value = 59
$ws1.Cells.Item(31, 7).Value = 75
$ws1.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82422"
$ws1.Cells.Item(31, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/ZzwSyv601709710699984.jpeg"

$ws1.Cells.Item(32, 1).Value = 31
$ws1.Cells.Item(32, 2).Value = "2024-04-13"
$ws1.Cells.Item(32, 3).Value = "杭州·代号鸢相聚广陵2.0only（中婚版）"
$ws1.Cells.Item(32, 4).Value = "凤起东路211号 名人名家宴会艺术中心(顺福店)"
$ws1.Cells.Item(32, 5).Value = "2024.04.13 10:00-04.13 17:00"
$ws1.Cells.Item(32, 6).Value = 133
$ws1.Cells.Item(32, 7).Value = 68
$ws1.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82176"
$ws1.Cells.Item(32, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/aWr1pXxU1709172854192.jpeg"

$ws1.Cells.Item(33, 1).Value = 32
$ws1.Cells.Item(33, 2).Value = "2024-04-13"
$ws1.Cells.Item(33, 3).Value = "杭州·赛马娘only—晴空雏菊"
$ws1.Cells.Item(33, 4).Value = "北干街道萧杭路689号 时尚外滩艺术中心"
$ws1.Cells.Item(33, 5).Value = "2024.04.13 09:00-04.13 18:00"
$ws1.Cells.Item(33, 6).Value = 199
$ws1.Cells.Item(33, 7).Value = 66
$ws1.Cells.Item(33, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81767"
$ws1.Cells.Item(33, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/ViMb8nbw1707122090281.jpeg"

$ws1.Cells.Item(34, 1).Value = 33
$ws1.Cells.Item(34, 2).Value = "2024-04-20"
$ws1.Cells.Item(34, 3).Value = "杭州·COMIC WORLD次元创作同人季特典·SP·浙里来消"
$ws1.Cells.Item(34, 4).Value = "德胜东路2539号 梦马汽车小镇"
$ws1.Cells.Item(34, 5).Value = "2024.04.20 10:00-04.21 17:00"
$ws1.Cells.Item(34, 6).Value = 18
$ws1.Cells.Item(34, 7).Value = 58
$ws1.Cells.Item(34, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82573"
$ws1.Cells.Item(34, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/UgyVGYJa1709879114323.png"

$ws1.Cells.Item(35, 1).Value = 34
$ws1.Cells.Item(35, 2).Value = "2024-04-20"
$ws1.Cells.Item(35, 3).Value = "杭州·SK怀旧展&偶像专场"
$ws1.Cells.Item(35, 4).Value = "沈半路171号 T-Car杭州汽车文化主题公园"
$ws1.Cells.Item(35, 5).Value = "2024.04.20 09:00-04.20 22:00"
$ws1.Cells.Item(35, 6).Value = 315
$ws1.Cells.Item(35, 7).Value = 60
$ws1.Cells.Item(35, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81764"
$ws1.Cells.Item(35, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/mtdbSuTZ1707119415384.jpeg"

$ws1.Cells.Item(36, 1).Value = 35
$ws1.Cells.Item(36, 2).Value = "2024-04-20"
$ws1.Cells.Item(36, 3).Value = "杭州·【海潮的回响Echo of The Tide】 | 刀客塔们的大群融入派对·明日方舟SPECIAL ONLY"
$ws1.Cells.Item(36, 4).Value = "保淑路2号 The Queen皇后"
$ws1.Cells.Item(36, 5).Value = "2024.04.20 14:00-04.20 18:00"
$ws1.Cells.Item(36, 6).Value = 61
$ws1.Cells.Item(36, 7).Value = 139
$ws1.Cells.Item(36, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82068"
$ws1.Cells.Item(36, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/gAR8Svsc1708919248994.png"

$ws1.Cells.Item(37, 1).Value = 36
$ws1.Cells.Item(37, 2).Value = "2024-04-20"
$ws1.Cells.Item(37, 3).Value = "杭州·白日梦次元动漫嘉年华"
$ws1.Cells.Item(37, 4).Value = "黄姑山路51-4号 0101park"
$ws1.Cells.Item(37, 5).Value = "2024.04.20 10:00-04.21 18:00"
$ws1.Cells.Item(37, 6).Value = 276
$ws1.Cells.Item(37, 7).Value = 68
$ws1.Cells.Item(37, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81634"
$ws1.Cells.Item(37, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/XnnWF6eP1709533743504.png"

$ws1.Cells.Item(38, 1).Value = 37
$ws1.Cells.Item(38, 2).Value = "2024-04-20"
$ws1.Cells.Item(38, 3).Value = "杭州·第五人格ONLY"
$ws1.Cells.Item(38, 4).Value = "望江东路333号 瑞莱克斯大酒店"
$ws1.Cells.Item(38, 5).Value = "2024.04.20 10:00-04.20 17:00"
$ws1.Cells.Item(38, 6).Value = 457
$ws1.Cells.Item(38, 7).Value = 60
$ws1.Cells.Item(38, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81987"
$ws1.Cells.Item(38, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/Dxk0hWDb1708572766103.jpeg"

$ws1.Cells.Item(39, 1).Value = 38
$ws1.Cells.Item(39, 2).Value = "2024-05-02"
$ws1.Cells.Item(39, 3).Value = "杭州·第四届华盟动漫次元嘉年华"
$ws1.Cells.Item(39, 4).Value = "创意路1号 中国智谷富春园区"
$ws1.Cells.Item(39, 5).Value = "2024.05.02 10:00-05.03 17:00"
$ws1.Cells.Item(39, 6).Value = 1316
$ws1.Cells.Item(39, 7).Value = 58
$ws1.Cells.Item(39, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82465"
$ws1.Cells.Item(39, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/4XHyqi3D1709780326858.jpeg"

$ws1.Cells.Item(40, 1).Value = 39
$ws1.Cells.Item(40, 2).Value = "2024-05-18"
$ws1.Cells.Item(40, 3).Value = "杭州·现世繁华-代号鸢only"
$ws1.Cells.Item(40, 4).Value = "石祥路575号 杭州海外海纳川大酒店(万达广场渡驾桥地铁站店)"
$ws1.Cells.Item(40, 5).Value = "2024.05.18 10:00-05.18 21:00"
$ws1.Cells.Item(40, 6).Value = 278
$ws1.Cells.Item(40, 7).Value = 76
$ws1.Cells.Item(40, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81905"
$ws1.Cells.Item(40, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/m3upuV2F1708327958926.jpeg"

$ws1.Cells.Item(41, 1).Value = 40
$ws1.Cells.Item(41, 2).Value = "2024-06-09"
$ws1.Cells.Item(41, 3).Value = "杭州·第三届日夜国乙only"
$ws1.Cells.Item(41, 4).Value = "创意路1号 中国智谷富春园区"
$ws1.Cells.Item(41, 5).Value = "2024.06.09 10:00-06.09 23:00"
$ws1.Cells.Item(41, 6).Value = 59
$ws1.Cells.Item(41, 7).Value = 58
$ws1.Cells.Item(41, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82618"
$ws1.Cells.Item(41, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/fXRzYEFH1710124366279.png"

$ws1.Cells.Item(42, 1).Value = 41
$ws1.Cells.Item(42, 2).Value = "2024-07-20"
$ws1.Cells.Item(42, 3).Value = "杭州·亚米二次茶话会展"
$ws1.Cells.Item(42, 4).Value = "湖州街20号 纳德自由酒店"
$ws1.Cells.Item(42, 5).Value = "2024.07.20 13:00-07.20 17:00"
$ws1.Cells.Item(42, 6).Value = 78
$ws1.Cells.Item(42, 7).Value = "不可售"
$ws1.Cells.Item(42, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81702"
$ws1.Cells.Item(42, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/uPDIsIoV1708311822716.jpeg"

$ws1.Cells.Item(43, 1).Value = 42
$ws1.Cells.Item(43, 2).Value = "2024-07-20"
$ws1.Cells.Item(43, 3).Value = "杭州·次元幻想--二次元全女夜场"
$ws1.Cells.Item(43, 4).Value = "保淑路2号 The Queen皇后"
$ws1.Cells.Item(43, 5).Value = "2024.07.20 13:00-07.20 19:00"
$ws1.Cells.Item(43, 6).Value = 279
$ws1.Cells.Item(43, 7).Value = 158
$ws1.Cells.Item(43, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81808"
$ws1.Cells.Item(43, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/sUUtSPh91707295826425.jpeg"

# ----- Sheet 2: "演出" (Performances) -----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(4, 6).Value = 12

# ----- Sheet 4: "全部类型" (All types) -----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 2592
$ws4.Cells.Item(5, 6).Value = 1449
$ws4.Cells.Item(7, 6).Value = 1124
$ws4.Cells.Item(11, 6).Value = 111
$ws4.Cells.Item(12, 6).Value = 549
$ws4.Cells.Item(13, 6).Value = 8934
$ws4.Cells.Item(14, 6).Value = 386
$ws4.Cells.Item(17, 6).Value = 253
$ws4.Cells.Item(22, 6).Value = 1166
$ws4.Cells.Item(24, 6).Value = 2058
$ws4.Cells.Item(25, 6).Value = 2135
$ws4.Cells.Item(27, 6).Value = 1839
$ws4.Cells.Item(31, 6).Value = 534
$ws4.Cells.Item(33, 6).Value = 133
$ws4.Cells.Item(34, 6).Value = 199
$ws4.Cells.Item(38, 6).Value = 276
$ws4.Cells.Item(39, 6).Value = 457
$ws4.Cells.Item(44, 6).Value = 1316
$ws4.Cells.Item(47, 6).Value = 59
$ws4.Cells.Item(49, 6).Value = 279
